$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for rows 2-9 from 45174 (2023-09-05)
# to 45175 (2023-09-06), keeping the existing date number format.
foreach ($row in 2..9) {
    $ws.Cells.Item($row, 3).Value = 45175
}
